$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Target end-of-document paragraph sequence (after this script runs):
#   15: Note that in the provided examples run times are printed. ...  (unchanged)
#   16: The default sparse=False argument in the reweighting methods ... (NEW)
#   17: (empty, ind left=360)                                           (unchanged)
#   18: (completely empty paragraph)                                   (NEW)
#   19: Iterative Enhanced Sampling (optional)                         (NEW position, no page break)
#   20: Relevant parameters: see Enhanced Sampling and Simulated Annealing (+ _GoBack bookmark) (NEW)
#   21: Instructions:  The enhanced sampling and simulated annealing ... (same text, re-flowed runs)
#   22: Conditional Probability Factorization (optional)                (NEW)
#   23: Relevant parameters: see Enhanced Sampling and Simulated Annealing (NEW)
#   24: Instructions: If one reaction coordinate is not sufficient ...   (NEW)
#   25: Note that the old probability used as an input for ...          (NEW)
# ---------------------------------------------------------------------------

# Step 1: insert the new "sparse=False" paragraph right after paragraph 15
#         ("Note that in the provided examples run times are printed...")
$pNote = $d.Paragraphs(15)
$pNote.Range.InsertParagraphAfter()
$pSparse = $d.Paragraphs(16)
$xmlSparse = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The default sparse=False argument in the reweighting methods will trim off negligible probabilities values from the ends of the returned probability array. </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>This argument is beneficial for SGOOP as it has trouble with sparse arrays. Sparse=True should be used when the probability is not being used as a SGOOP input (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> input for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rebiasing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> or counting wells). Sparse=True may also be beneficial for systems that have important features with very low probability as these would be neglected otherwise.</w:t></w:r></w:p>
'@
$pSparse.Range.InsertXML($xmlSparse)

# Step 2: remove the two now-obsolete paragraphs that used to follow the blank
#         (ind=360) paragraph: the bookmark-only paragraph and the
#         "Iterative Enhanced Sampling (optional)" heading paragraph. Their
#         content will be rebuilt fresh, in the correct new order, below.
#         After step 1, paragraph numbering is:
#           16 = sparse (new)      17 = blank ind=360 (old 16)
#           18 = bookmark-only (old 17)   19 = Iterative Enhanced Sampling (old 18)
#           20 = Instructions: enhanced sampling... (old 19)
$d.Paragraphs(18).Range.Delete()
$d.Paragraphs(18).Range.Delete()
# Now: 16 = sparse, 17 = blank ind=360, 18 = Instructions: enhanced sampling...(old 19)

# Step 3: insert three fresh paragraphs between the blank (ind=360) paragraph
#         and the Instructions paragraph: a completely empty paragraph, the
#         "Iterative Enhanced Sampling (optional)" heading, and the
#         "Relevant parameters..." bullet (carrying the _GoBack bookmark).
$pBlank = $d.Paragraphs(17)
$pBlank.Range.InsertParagraphAfter()
$pBlank.Range.InsertParagraphAfter()
$pBlank.Range.InsertParagraphAfter()

$pEmpty = $d.Paragraphs(18)
$xmlEmpty = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$pEmpty.Range.InsertXML($xmlEmpty)

$pIterHeading = $d.Paragraphs(19)
$xmlIterHeading = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Iterative Enhanced Sampling (optional)</w:t></w:r></w:p>
'@
$pIterHeading.Range.InsertXML($xmlIterHeading)

$pRelevant = $d.Paragraphs(20)
$xmlRelevant = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Relevant parameters: see Enhanced Sampling and Simulated Annealing</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$pRelevant.Range.InsertXML($xmlRelevant)

# Step 4: rewrite the "Instructions" paragraph (now at index 21) so its runs
#         match the target re-flow (same visible text, merged trailing runs).
$pInstructions = $d.Paragraphs(21)
$xmlInstructions = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Instructions:  The enhanced sampling and simulated annealing steps can be repeated using the SGOOP reaction coordinate to add bias.</w:t></w:r><w:r><w:t xml:space="preserve">  Change the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>metadynamics</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> reaction coordinate in PLUMED to match the SGOOP reaction coordinate.  Repeat the simulated annealing process with the new COLVAR and FES files from subsequent </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>metadynamics</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> runs while using the same unbiased COLVAR file for maximum caliber calculations.  </w:t></w:r><w:r><w:t>This will improve sampling allow for more accurate SGOOP calculations.  After a few iterations of SGOOP and molecular dynamics, the reaction coordinate should converge.</w:t></w:r></w:p>
'@
$pInstructions.Range.InsertXML($xmlInstructions)

# Step 5: append the four brand-new paragraphs describing "Conditional
#         Probability Factorization (optional)" right after the Instructions
#         paragraph.
$pInstructions.Range.InsertParagraphAfter()
$pInstructions.Range.InsertParagraphAfter()
$pInstructions.Range.InsertParagraphAfter()
$pInstructions.Range.InsertParagraphAfter()

$pCondHeading = $d.Paragraphs(22)
$xmlCondHeading = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Conditional Probability Factorization (optional)</w:t></w:r></w:p>
'@
$pCondHeading.Range.InsertXML($xmlCondHeading)

$pCondRelevant = $d.Paragraphs(23)
$xmlCondRelevant = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Relevant parameters: </w:t></w:r><w:r><w:t>see Enhanced Sampling and Simulated Annealing</w:t></w:r></w:p>
'@
$pCondRelevant.Range.InsertXML($xmlCondRelevant)

$pCondInstructions = $d.Paragraphs(24)
$xmlCondInstructions = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Instructions: I</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>f one reaction coordinate is not sufficient for sampling (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> some features are hidden) more reaction coordinates may be generated. Use the function </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rebias</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> instead of the function reweight to generate a probability conditional on the first reaction coordinate. A trajectory with bias along the first reaction coordinate should be used for maximum caliber as this trajectory will have a free energy surface that is close to that of the conditional probability.</w:t></w:r></w:p>
'@
$pCondInstructions.Range.InsertXML($xmlCondInstructions)

$pCondNote = $d.Paragraphs(25)
$xmlCondNote = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Note that the old probability used as an input for </w:t></w:r><w:r><w:t>the conditional probability should be generated using the sparse=True argument in order to account for the full probability along the first reaction coordinate.</w:t></w:r></w:p>
'@
$pCondNote.Range.InsertXML($xmlCondNote)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
